$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.899.39'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '1.584.45'
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.17'
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("E7").Value = '  -2.04%  '
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("E9").Value = '  -1.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.08'
$ws.Range("E10").Value = '  -0.40%  '
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("E12").Value = '  -1.86%  '
$ws.Range("D13").Value = '1.580.76'
$ws.Range("E13").Value = '  -2.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.02'
$ws.Range("E14").Value = '  -1.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.506'
$ws.Range("E15").Value = '  -1.97%  '
$ws.Range("D16").Value = '25.889.27'
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D17").Value = '0.0₃0725'
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '59.96'
$ws.Range("E18").Value = '  -2.43%  '
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '192.97'
$ws.Range("E20").Value = '  +1.03%  '
$ws.Range("E21").Value = '  -0.40%  '
$ws.Range("E22").Value = '  -0.54%  '
$ws.Range("E23").Value = '  -0.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.130'
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.53'
$ws.Range("E25").Value = '  -1.27%  '
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.44'
$ws.Range("E29").Value = '  -2.25%  '
$ws.Range("E30").Value = '  -4.71%  '
$ws.Range("E31").Value = '  -0.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.11'
$ws.Range("E32").Value = '  +0.46%  '
$ws.Range("E33").Value = '  -1.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.49'
$ws.Range("E34").Value = '  +1.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.36'
$ws.Range("E35").Value = '  -2.22%  '
$ws.Range("D36").Value = '1.096.17'
$ws.Range("E36").Value = '  -2.38%  '
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("E38").Value = '  -2.09%  '
$ws.Range("E39").Value = '  -0.23%  '
$ws.Range("E40").Value = '  -2.29%  '
$ws.Range("E41").Value = '  -4.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.800'
$ws.Range("E42").Value = '  +5.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '93.46'
$ws.Range("E43").Value = '  -3.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.11'
$ws.Range("E44").Value = '  +0.84%  '
$ws.Range("D45").Value = '1.719.07'
$ws.Range("E45").Value = '  -1.86%  '
$ws.Range("E46").Value = '  -1.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.50'
$ws.Range("E47").Value = '  +2.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.13'
$ws.Range("E48").Value = '  -0.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0509'
$ws.Range("E49").Value = '  -1.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.407'
$ws.Range("E50").Value = '  -1.04%  '
$ws.Range("E51").Value = '  -0.22%  '
